# Combine dc parameters into single file
# Adds a "Distance / Distance^2 / Distance^3" coefficient block (rows 6-8)
# to the existing HBW/HBO/NHB parameter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pre-stage the new cell style (font 1, centered, wrap text) on a ---
# --- scratch cell so it gets allocated before the temporary "Text"   ---
# --- number-format style used below (keeps the style table tidy).   ---
$ws.Range("B1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").Clear()

# --- The regression coefficients must stay as literal text (they carry ---
# --- trailing spaces, e.g. "-0.0801 "), so force the cells to Text    ---
# --- format before writing them - otherwise Excel auto-converts them  ---
# --- to numbers and the trailing space is lost.                      ---
$ws.Range("B6:E8").NumberFormat = "@"

$ws.Range("B6").Value = "-0.0801 "
$ws.Range("C6").Value = "-0.1728 "
$ws.Range("D6").Value = "-0.1157 "
$ws.Range("E6").Value = "0.3151 "

$ws.Range("B7").Value = "0.0026 "
$ws.Range("C7").Value = "0.0034 "
$ws.Range("D7").Value = "0.0035 "
$ws.Range("E7").Value = "-0.0026 "

$ws.Range("B8").Value = "-0.0000090 "
$ws.Range("C8").Value = "-0.0000110 "
$ws.Range("D8").Value = "-0.0000133 "
$ws.Range("E8").Value = "0.0000055 "

$ws.Range("A6").Value = "Distance"
$ws.Range("A7").Value = "Distance^2"
$ws.Range("A8").Value = "Distance^3"

# --- Apply the new row's look: same font as the rest of the table, ---
# --- centered horizontally/vertically, with wrapped text.          ---
$ws.Range("B1").Copy()
$ws.Range("A6:E8").PasteSpecial(-4122)
$ws.Range("A6:E8").WrapText = $true

# --- Match the taller header-style rows used for this new block. ---
$ws.Rows("6:8").RowHeight = 17

# --- Leave the selection where Excel would land after typing the block. ---
$null = $ws.Range("A9").Select()
